$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (A = date serial, C = indicatore_stress_t, D = indicatore_stress_t1)
# Rows 8 and 14 lose their C value (cell becomes empty) while D takes the old C value.
$rows = @(
    @{ Row = 2;  A = 44409; C = 0.37;               D = 0.39 },
    @{ Row = 3;  A = 44409; C = 0.39;               D = 0.36 },
    @{ Row = 4;  A = 44409; C = 0.49;               D = 0.47 },
    @{ Row = 5;  A = 44409; C = 0.37;               D = 0.36 },
    @{ Row = 6;  A = 44409; C = 0.77;               D = 0.77 },
    @{ Row = 7;  A = 44409; C = 0.41;               D = 0.4  },
    @{ Row = 8;  A = 44409; C = $null;              D = 0.74 },
    @{ Row = 9;  A = 44409; C = 0.61;               D = 0.5600000000000001 },
    @{ Row = 10; A = 44409; C = 0.38;               D = 0.38 },
    @{ Row = 11; A = 44409; C = 0.52;               D = 0.5  },
    @{ Row = 12; A = 44409; C = 0.19;               D = 0.17 },
    @{ Row = 13; A = 44409; C = 0.33;               D = 0.34 },
    @{ Row = 14; A = 44409; C = $null;              D = 0.45 },
    @{ Row = 15; A = 44409; C = 0.31;               D = 0.3  },
    @{ Row = 16; A = 44409; C = 0.27;               D = 0.27 },
    @{ Row = 17; A = 44409; C = 1.39;               D = 1.39 },
    @{ Row = 18; A = 44409; C = 0.97;               D = 0.96 },
    @{ Row = 19; A = 44409; C = 1.17;               D = 1.11 },
    @{ Row = 20; A = 44409; C = 0.83;               D = 0.86 },
    @{ Row = 21; A = 44409; C = 0.22;               D = 0.23 },
    @{ Row = 22; A = 44409; C = 0.93;               D = 0.9399999999999999 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    if ($null -eq $r.C) {
        $ws.Cells.Item($r.Row, 3).ClearContents()
    } else {
        $ws.Cells.Item($r.Row, 3).Value = $r.C
    }
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
